$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "nieve"
$ws.Range("B2").Value = 36
$ws.Range("C2").Value = 100
